$wb = $excel.ActiveWorkbook

# --- Sheet "Veicolo": update id in A2 ---
$wsVeicolo = $wb.Worksheets.Item("Veicolo")
$wsVeicolo.Range("A2").Value = 253940

# --- Sheet "Release Date (RD)": update rows 2-3, remove rows 4-17 ---
$wsRD = $wb.Worksheets.Item("Release Date (RD)")

# Delete old rows 4 through 17 (16 rows), leaving header + 2 data rows
$wsRD.Range("A4:C17").EntireRow.Delete() | Out-Null

# Update remaining data rows 2 and 3
$wsRD.Range("A2").Value = 254187
$wsRD.Range("B2").Value = 45975.58333333334
$wsRD.Range("C2").Value = 0

$wsRD.Range("A3").Value = 254967
$wsRD.Range("B3").Value = 45981.58333333334
$wsRD.Range("C3").Value = 0

# --- Sheet "RD Tassative": update row 2, add rows 3-4 ---
$wsTass = $wb.Worksheets.Item("RD Tassative")

$wsTass.Range("A2").Value = 254237
$wsTass.Range("B2").Value = 45975.58333333334

$wsTass.Range("A3").Value = 254427
$wsTass.Range("B3").Value = 45978.58333333334
$wsTass.Range("B3").NumberFormat = $wsTass.Range("B2").NumberFormat

$wsTass.Range("A4").Value = 254428
$wsTass.Range("B4").Value = 45978.58333333334
$wsTass.Range("B4").NumberFormat = $wsTass.Range("B2").NumberFormat
